# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 75
$ws1.Range("F3").Value = 804
$ws1.Range("F6").Value = 108
$ws1.Range("F8").Value = 4355
$ws1.Range("F9").Value = 93
$ws1.Range("F10").Value = 4938
$ws1.Range("F11").Value = 549
$ws1.Range("F12").Value = 1242
$ws1.Range("F13").Value = 84

# Sheet "全部类型" (sheetId 4) - same events, but row numbers are shifted
# by one extra row (a 演出 event inserted at row 7) relative to 展览
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 75
$ws4.Range("F3").Value = 804
$ws4.Range("F6").Value = 108
$ws4.Range("F9").Value = 4355
$ws4.Range("F10").Value = 93
$ws4.Range("F11").Value = 4938
$ws4.Range("F12").Value = 549
$ws4.Range("F13").Value = 1242
$ws4.Range("F14").Value = 84
